# Fix saleTable rows 12-36: correct singleCost (barcode/name/price) mismatches
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "'060143"
$ws.Range("B12").Value = "苏泊尔水煲"
$ws.Range("D12").Value = 119

$ws.Range("A13").Value = "'090150"
$ws.Range("B13").Value = "水煲"
$ws.Range("D13").Value = 99

$ws.Range("A14").Value = "'090163"
$ws.Range("B14").Value = "永兴3L火锅JQH-100"
$ws.Range("D14").Value = 179

$ws.Range("A15").Value = "'090164"
$ws.Range("B15").Value = "三角电钣饭锅"
$ws.Range("D15").Value = 98

$ws.Range("A16").Value = "'090359"
$ws.Range("B16").Value = "半球电水煲2L"
$ws.Range("D16").Value = 39

$ws.Range("A17").Value = "'135233234"
$ws.Range("B17").Value = "闹钟"
$ws.Range("D17").Value = 18

$ws.Range("A18").Value = "'135015"
$ws.Range("B18").Value = "LED彩色小夜灯"
$ws.Range("D18").Value = 5

$ws.Range("A19").Value = "'6927065410459"
$ws.Range("B19").Value = "1021空调扇"
$ws.Range("D19").Value = 350

$ws.Range("A20").Value = "'135012"
$ws.Range("B20").Value = "低碳无福射节能灯（大）"
$ws.Range("D20").Value = 12

$ws.Range("A21").Value = "'135011"
$ws.Range("B21").Value = "万意款节王4#节能双灶"
$ws.Range("D21").Value = 78.4

$ws.Range("A22").Value = "'6927065400788"
$ws.Range("B22").Value = "DL0740电风扇"
$ws.Range("D22").Value = 30

$ws.Range("A23").Value = "'6926159300034"
$ws.Range("B23").Value = "低碳节能小夜灯"
$ws.Range("D23").Value = 9.9

$ws.Range("A24").Value = "'6950610208522"
$ws.Range("B24").Value = "苏泊尔电水煲1702A"
$ws.Range("D24").Value = 169

$ws.Range("A25").Value = "'041184"
$ws.Range("B25").Value = "炊大皇不沾锅"
$ws.Range("D25").Value = 188

$ws.Range("A26").Value = "'060151"
$ws.Range("B26").Value = "苏泊尔电饼档"
$ws.Range("D26").Value = 229

$ws.Range("A27").Value = "'060158"
$ws.Range("B27").Value = "苏泊尔电饭煲"
$ws.Range("D27").Value = 159

$ws.Range("A28").Value = "'060168"
$ws.Range("B28").Value = "苏泊尔迷你电饭煲"
$ws.Range("D28").Value = 199

$ws.Range("A29").Value = "'090149"
$ws.Range("B29").Value = "电压力锅"
$ws.Range("D29").Value = 499

$ws.Range("A30").Value = "'090165"
$ws.Range("B30").Value = "三角电饭锅"
$ws.Range("D30").Value = 118

$ws.Range("A31").Value = "'090166"
$ws.Range("B31").Value = "三角电饭锅"
$ws.Range("D31").Value = 138

$ws.Range("A32").Value = "'091456"
$ws.Range("B32").Value = "FSJ-180学生扇"
$ws.Range("D32").Value = 57.8

$ws.Range("A33").Value = "'091469"
$ws.Range("B33").Value = "1535电风扇"
$ws.Range("D33").Value = 169

$ws.Range("A34").Value = "'091514"
$ws.Range("B34").Value = "低碳无福射节能灯（大）"
$ws.Range("D34").Value = 12

$ws.Range("A35").Value = "'091545"
$ws.Range("B35").Value = "LED彩色小夜灯"
$ws.Range("D35").Value = 5

$ws.Range("A36").Value = "'091546"
$ws.Range("B36").Value = "低碳节能小夜灯"
$ws.Range("D36").Value = 8
